$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 / J1, matching the formatting of the existing header (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New data values for columns I (I0) and J (IF)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 5

$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 7

$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 9
